# Template_Questions.xlsx update:
#  - fix wording in two shared strings on the "Questions" sheet
#  - widen column B on the "Textes" sheet
#  - add a new "Introduction" sheet (RGPD / presentation text) and make it active

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Questions" sheet: correct wording ("un champ \"Autre\"" -> "une reponse
#    autre") in the two explanatory header cells.
# ---------------------------------------------------------------------------
$wsQuestions = $wb.Worksheets.Item("Questions")

$wsQuestions.Range("B1").Value = 'Type de réponse attendu : "text", "radio" (choix unique), "checkbox" (choix multiple). Si vous souhaitez indiquer une réponse autre à votre réponse à choix multiple ou unique, indiquer : "checkbox, text" ou "radio, text"'
$wsQuestions.Range("D1").Value = 'Créer autant de colonnes que de réponses possible (avec leurs noms). Attention si vous avez une réponse autre ne l''ajouté pas en temps que colonne.'

$wsQuestions.Range("C4").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) "Textes" sheet: widen column B so the longer labels fit.
# ---------------------------------------------------------------------------
$wsTextes = $wb.Worksheets.Item("Textes")
$wsTextes.Columns.Item(2).ColumnWidth = 78.5

$wsTextes.Range("C2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Add the new "Introduction" sheet after "Textes".
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsIntro = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsIntro.Name = "Introduction"

$wsIntro.Columns.Item(1).ColumnWidth = 26.27
$wsIntro.Columns.Item(2).ColumnWidth = 82.8

# Row 1: RGPD
$wsIntro.Rows.Item(1).RowHeight = 129
$a1 = $wsIntro.Range("A1")
$a1.Value = "RGPD"
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4108
$a1.WrapText = $true
$a1.Font.Color = 204

$b1 = $wsIntro.Range("B1")
$b1.Value = "Ce questionnaire effectue une récolte de vos données personnelles (caractéristiques physiques, goûts, coordonnées notamment) ainsi que de données sur la manière dont vous le remplissez (vitesse, hésitations, changements, etc.). Ces données seront utilisées à des fins de recherches dans le cadre des travaux de Mme Jessica Pidoux. Pour toute demande concernant la gestion des données personnelles, envoyer un mail à l'adresse 'jessica.pidoux@epfl.ch'"
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4108
$b1.WrapText = $true

# Row 2: Validation du RGPD
$wsIntro.Rows.Item(2).RowHeight = 87.9
$a2 = $wsIntro.Range("A2")
$a2.Value = "Validation du RGPD"
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4108
$a2.WrapText = $true
$a2.Font.Color = 204

$b2 = $wsIntro.Range("B2")
$b2.Value = "Cochez la case si vous acceptez les conditions ci-dessus"
$b2.HorizontalAlignment = -4108
$b2.VerticalAlignment = -4108

# Row 3: Presentation
$wsIntro.Rows.Item(3).RowHeight = 97.5
$a3 = $wsIntro.Range("A3")
$a3.Value = "Presentation"
$a3.HorizontalAlignment = -4108
$a3.VerticalAlignment = -4108
$a3.WrapText = $true
$a3.Font.Color = 204

$b3 = $wsIntro.Range("B3")
$b3.Value = "Le questionnaire est divisé en plusieurs phases : une phase de questions de départ sur votre sexe et votre orientation sexuelle, plusieurs phases de classement, une phase de questions supplémentaires à la fin. Lors des phases de classement, vous devrez déplacer avec la souris des descriptions de caractéristiques vers une échelle d'importance."
$b3.HorizontalAlignment = -4108
$b3.VerticalAlignment = -4108
$b3.WrapText = $true

$wsIntro.Activate() | Out-Null
$wsIntro.Range("B3").Select() | Out-Null
